$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "62.174.57"
$ws.Range("E2").Value = "  +2.85%  "

# Row 3
$ws.Range("D3").Value = "2.420.05"
$ws.Range("E3").Value = "  +3.64%  "

# Row 4
$ws.Range("E4").Value = "  +0.20%  "

# Row 5
$ws.Range("D5").Value = "'554.88"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +2.45%  "

# Row 6
$ws.Range("D6").Value = "'143.24"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +5.75%  "

# Row 7
$ws.Range("E7").Value = "  +0.21%  "

# Row 8
$ws.Range("D8").Value = "'0.533"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +2.70%  "

# Row 9
$ws.Range("D9").Value = "2.419.21"
$ws.Range("E9").Value = "  +3.67%  "

# Row 10
$ws.Range("E10").Value = "  +5.61%  "

# Row 11
$ws.Range("E11").Value = "  +1.11%  "

# Row 12
$ws.Range("D12").Value = "'5.40"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +2.66%  "

# Row 13
$ws.Range("D13").Value = "'0.353"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +4.84%  "

# Row 14
$ws.Range("D14").Value = "'26.30"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +8.14%  "

# Row 15
$ws.Range("E15").Value = "  +10.59%  "

# Row 16
$ws.Range("D16").Value = "2.858.05"
$ws.Range("E16").Value = "  +3.98%  "

# Row 17
$ws.Range("D17").Value = "61.982.41"
$ws.Range("E17").Value = "  +3.03%  "

# Row 18
$ws.Range("D18").Value = "2.419.67"
$ws.Range("E18").Value = "  +3.88%  "

# Row 19
$ws.Range("D19").Value = "'11.14"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +5.94%  "

# Row 20
$ws.Range("D20").Value = "'324.87"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +2.84%  "

# Row 21
$ws.Range("D21").Value = "'4.19"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +3.48%  "

# Row 22
$ws.Range("D22").Value = "'6.75"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +4.14%  "

# Row 23
$ws.Range("E23").Value = "  +0.23%  "

# Row 24
$ws.Range("D24").Value = "'64.91"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +3.67%  "

# Row 25
$ws.Range("E25").Value = "  +6.30%  "

# Row 26
$ws.Range("D26").Value = "'9.22"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +11.45%  "

# Row 27
$ws.Range("D27").Value = "'565.03"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +15.60%  "

# Row 28
$ws.Range("B28").Value = "WrappedeETH"
$ws.Range("C28").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D28").Value = "2.541.85"
$ws.Range("E28").Value = "  +4.18%  "

# Row 29
$ws.Range("B29").Value = "Binance-PegBSC-USD"
$ws.Range("C29").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D29").Value = "'0.999"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.11%  "

# Row 30
$ws.Range("B30").Value = "PEPE"
$ws.Range("C30").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D30").Value = "0.0₃0939"
$ws.Range("E30").Value = "  +10.78%  "

# Row 31
$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").Value = "'8.36"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +6.83%  "

# Row 32
$ws.Range("E32").Value = "  +7.41%  "

# Row 33
$ws.Range("D33").Value = "'0.149"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +3.86%  "

# Row 34
$ws.Range("D34").Value = "'1.87"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +5.30%  "

# Row 35
$ws.Range("D35").Value = "'1.56"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +4.95%  "

# Row 37
$ws.Range("E37").Value = "  +12.85%  "

# Row 38
$ws.Range("E38").Value = "  +0.39%  "

# Row 39
$ws.Range("D39").Value = "'4.82"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +6.87%  "

# Row 40
$ws.Range("E40").Value = "  +3.55%  "

# Row 41
$ws.Range("D41").Value = "'18.81"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +2.27%  "

# Row 42
$ws.Range("D42").Value = "'147.15"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +4.16%  "

# Row 43
$ws.Range("E43").Value = "  +0.32%  "

# Row 44
$ws.Range("D44").Value = "'2.32"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +15.01%  "

# Row 45
$ws.Range("D45").Value = "'151.57"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +8.13%  "

# Row 46
$ws.Range("D46").Value = "'3.64"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +3.39%  "

# Row 47
$ws.Range("D47").Value = "'0.0540"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +6.65%  "

# Row 48
$ws.Range("D48").Value = "'20.39"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +8.05%  "

# Row 49
$ws.Range("D49").Value = "'0.590"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +4.66%  "

# Row 50
$ws.Range("B50").Value = "Stellar"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D50").Value = "'0.0910"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +2.06%  "

# Row 51
$ws.Range("B51").Value = "VeChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D51").Value = "'0.0227"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +4.44%  "
